$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - update UCRN + INSIGHT REF + MATTER TYPE (ASSIGNED BY / TASK stay the same)
$ws.Range("C2").Value = "S459/50672"
$ws.Range("D2").Value = "NNT1028"
$ws.Range("E2").Value = "New Build Purchase"

# Row 3 - new row
$ws.Range("A3").Value = "AU"
$ws.Range("B3").Value = "SDLT Submitted"
$ws.Range("C3").Value = "S459/53841"
$ws.Range("D3").Value = "NBT1872"
$ws.Range("E3").Value = "New Build Purchase"

# Row 4 - new row
$ws.Range("A4").Value = "AU"
$ws.Range("B4").Value = "SDLT Submitted"
$ws.Range("C4").Value = "S459/53760"
$ws.Range("D4").Value = "NNT1227"
$ws.Range("E4").Value = "New Build Purchase"

# Update selection to match the diff's recorded view state
$ws.Range("D11").Select()
